# "Generate Report for Archive"
#
# 1) Status text: every "Ready for handoff" cell (Overview!E2:F3, zh-cn!C2:C3,
#    de-de!C2:C3 all share one string) becomes "In Translation".
# 2) The now-shorter status text no longer needs as wide a column, so the
#    "Status" columns are narrowed:
#      - Overview: columns E (zh-cn) and F (de-de)
#      - zh-cn:    column C (Status)
#      - de-de:    column C (Status)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
